$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.07901772845256971
$ws.Range("D2").Value = 0.01282090158601434
$ws.Range("E2").Value = 0.4189125050028082
$ws.Range("F2").Value = 0.7907120707850055
$ws.Range("G2").Value = 0.6460626888141263
$ws.Range("H2").Value = 0.697944098715368
$ws.Range("I2").Value = 0.4417292335453524
$ws.Range("K2").Value = 0.9892845052421819
$ws.Range("B3").Value = 0.06996419166664225
$ws.Range("D3").Value = 0.01276061562164088
$ws.Range("E3").Value = 0.3652792235888711
$ws.Range("F3").Value = 0.7733542935691702
$ws.Range("G3").Value = 0.6281733772643889
$ws.Range("H3").Value = 0.6956424957919722
$ws.Range("I3").Value = 0.4493519771821184
$ws.Range("K3").Value = 0.863236354148313
$ws.Range("B4").Value = 0.06439910882150457
$ws.Range("D4").Value = 0.01272669363102708
$ws.Range("E4").Value = 0.3324627557361595
$ws.Range("F4").Value = 0.7633976108373872
$ws.Range("G4").Value = 0.6178401441101045
$ws.Range("H4").Value = 0.6947415389205105
$ws.Range("I4").Value = 0.4543193649620427
$ws.Range("K4").Value = 0.7857583879169283
$ws.Range("B5").Value = 0.06212994706659458
$ws.Range("D5").Value = 0.01271364995692181
$ws.Range("E5").Value = 0.3191161875266459
$ws.Range("F5").Value = 0.7595152273663075
$ws.Range("G5").Value = 0.6137914163646201
$ws.Range("H5").Value = 0.6945026234747615
$ws.Range("I5").Value = 0.456415642640005
$ws.Range("K5").Value = 0.7541639360959209
$ws.Range("B6").Value = 0.06175307937374441
$ws.Range("D6").Value = 0.01271153122995727
$ws.Range("E6").Value = 0.3169015252625229
$ws.Range("F6").Value = 0.7588810966795592
$ws.Range("G6").Value = 0.6131288772888581
$ws.Range("H6").Value = 0.6944706799173872
$ws.Range("I6").Value = 0.4567680741067441
$ws.Range("K6").Value = 0.7489163872529048
$ws.Range("B7").Value = 0.06436851128292176
$ws.Range("D7").Value = 0.01272651456006813
$ws.Range("E7").Value = 0.3322826554600766
$ws.Range("F7").Value = 0.763344544552865
$ws.Range("G7").Value = 0.6177848870316325
$ws.Range("H7").Value = 0.6947377983867398
$ws.Range("I7").Value = 0.4543473446507456
$ws.Range("K7").Value = 0.7853323819221885
$ws.Range("B8").Value = 0.07589747143293835
$ws.Range("D8").Value = 0.01279947372318802
$ws.Range("E8").Value = 0.4003944994745865
$ws.Range("F8").Value = 0.7845808465062873
$ws.Range("G8").Value = 0.6397584411409412
$ws.Range("H8").Value = 0.6970438416886253
$ws.Range("I8").Value = 0.444297907042424
$ws.Range("K8").Value = 0.9458401928116018
$ws.Range("B9").Value = 0.09844880170884096
$ws.Range("D9").Value = 0.012967064628576
$ws.Range("E9").Value = 0.5349859246698117
$ws.Range("F9").Value = 0.8318472573572251
$ws.Range("G9").Value = 0.6880846735884631
$ws.Range("H9").Value = 0.7056576607080558
$ws.Range("I9").Value = 0.4268741312777813
$ws.Range("K9").Value = 1.259967482288346
$ws.Range("B10").Value = 0.1149736298370527
$ws.Range("D10").Value = 0.01310517650616561
$ws.Range("E10").Value = 0.6346638144554788
$ws.Range("F10").Value = 0.8700885669360616
$ws.Range("G10").Value = 0.7268860319001362
$ws.Range("H10").Value = 0.7145204741756856
$ws.Range("I10").Value = 0.4154725780895436
$ws.Range("K10").Value = 1.490450583387201
$ws.Range("B11").Value = 0.1224799730716342
$ws.Range("D11").Value = 0.01317128528015843
$ws.Range("E11").Value = 0.6802196872274351
$ws.Range("F11").Value = 0.8882679699421487
$ws.Range("G11").Value = 0.7452763281368107
$ws.Range("H11").Value = 0.7191114946402308
$ws.Range("I11").Value = 0.4105913267349006
$ws.Range("K11").Value = 1.59525535393476
$ws.Range("B12").Value = 0.1253206852089193
$ws.Range("D12").Value = 0.01319679336254254
$ws.Range("E12").Value = 0.6975038901726975
$ws.Range("F12").Value = 0.8952661498251473
$ws.Range("G12").Value = 0.7523483567790947
$ws.Range("H12").Value = 0.7209311110966894
$ws.Range("I12").Value = 0.4087869892821416
$ws.Range("K12").Value = 1.634937034750465
$ws.Range("B13").Value = 0.1247089694124952
$ws.Range("D13").Value = 0.01319127860814007
$ws.Range("E13").Value = 0.6937799083939922
$ws.Range("F13").Value = 0.8937538723284746
$ws.Range("G13").Value = 0.7508204372615523
$ws.Range("H13").Value = 0.7205356065443027
$ws.Range("I13").Value = 0.409173622090135
$ws.Range("K13").Value = 1.626391114481294
$ws.Range("B14").Value = 0.1227137169917114
$ws.Range("D14").Value = 0.01317337432727683
$ws.Range("E14").Value = 0.6816409882942622
$ws.Range("F14").Value = 0.8888414208924473
$ws.Range("G14").Value = 0.7458559734475614
$ws.Range("H14").Value = 0.7192595664041903
$ws.Range("I14").Value = 0.4104419982458403
$ws.Range("K14").Value = 1.598520098614188
$ws.Range("B15").Value = 0.1214913295848135
$ws.Range("D15").Value = 0.01316246926486286
$ws.Range("E15").Value = 0.6742099517796447
$ws.Range("F15").Value = 0.8858472951498584
$ws.Range("G15").Value = 0.7428292145440594
$ws.Range("H15").Value = 0.7184885362422051
$ws.Range("I15").Value = 0.4112246618182915
$ws.Range("K15").Value = 1.581447584178875
$ws.Range("B16").Value = 0.1144828354586451
$ws.Range("D16").Value = 0.01310092234574967
$ws.Range("E16").Value = 0.6316911258656717
$ws.Range("F16").Value = 0.8689163808269882
$ws.Range("G16").Value = 0.7256992059120364
$ws.Range("H16").Value = 0.7142317516916137
$ws.Range("I16").Value = 0.4157977390970728
$ws.Range("K16").Value = 1.483600526730868
$ws.Range("B17").Value = 0.1101804120322214
$ws.Range("D17").Value = 0.01306400705365718
$ws.Range("E17").Value = 0.6056633925233825
$ws.Range("F17").Value = 0.8587314155990242
$ws.Range("G17").Value = 0.7153810762138448
$ws.Range("H17").Value = 0.7117640739007811
$ws.Range("I17").Value = 0.418681510792851
$ws.Range("K17").Value = 1.423563884767191
$ws.Range("B18").Value = 0.1077047599404324
$ws.Range("D18").Value = 0.01304308298910328
$ws.Range("E18").Value = 0.590712709172351
$ws.Range("F18").Value = 0.8529469303609858
$ws.Range("G18").Value = 0.709515850626957
$ws.Range("H18").Value = 0.7103973325230584
$ws.Range("I18").Value = 0.4203689126210719
$ws.Range("K18").Value = 1.389028296960191
$ws.Range("B19").Value = 0.1068663793808895
$ws.Range("D19").Value = 0.01303605141087871
$ws.Range("E19").Value = 0.5856539800581402
$ws.Range("F19").Value = 0.851001007721365
$ws.Range("G19").Value = 0.7075418692942037
$ws.Range("H19").Value = 0.7099435913806929
$ws.Range("I19").Value = 0.420945167442234
$ws.Range("K19").Value = 1.377334418853252
$ws.Range("B20").Value = 0.1106385186241567
$ws.Range("D20").Value = 0.01306790479326381
$ws.Range("E20").Value = 0.608432022405097
$ws.Range("F20").Value = 0.8598079912279672
$ws.Range("G20").Value = 0.7164722545417419
$ws.Range("H20").Value = 0.7120213135237918
$ws.Range("I20").Value = 0.4183715534062511
$ws.Range("K20").Value = 1.42995530785123
$ws.Range("B21").Value = 0.1232998209021474
$ws.Range("D21").Value = 0.01317862035873318
$ws.Range("E21").Value = 0.685205560997133
$ws.Range("F21").Value = 0.8902812208248889
$ws.Range("G21").Value = 0.7473112114606408
$ws.Range("H21").Value = 0.7196321639859775
$ws.Range("I21").Value = 0.4100682475336654
$ws.Range("K21").Value = 1.606706638939499
$ws.Range("B22").Value = 0.1315642586652501
$ws.Range("D22").Value = 0.01325374456342487
$ws.Range("E22").Value = 0.7355762087192659
$ws.Range("F22").Value = 0.9108625804280166
$ws.Range("G22").Value = 0.7680966920417518
$ws.Range("H22").Value = 0.7250793057865508
$ws.Range("I22").Value = 0.4048986086142445
$ws.Range("K22").Value = 1.722191448607475
$ws.Range("B23").Value = 0.1271543987773356
$ws.Range("D23").Value = 0.01321339538021604
$ws.Range("E23").Value = 0.7086737239511223
$ws.Range("F23").Value = 0.8998165745952633
$ws.Range("G23").Value = 0.7569448379554728
$ws.Range("H23").Value = 0.7221285588732655
$ws.Range("I23").Value = 0.4076341641078489
$ws.Range("K23").Value = 1.660557790748101
$ws.Range("B24").Value = 0.110431415257878
$ws.Range("D24").Value = 0.01306614169263298
$ws.Range("E24").Value = 0.6071802845537633
$ws.Range("F24").Value = 0.8593210503729978
$ws.Range("G24").Value = 0.7159787248587577
$ws.Range("H24").Value = 0.7119048536863488
$ws.Range("I24").Value = 0.4185115933157562
$ws.Range("K24").Value = 1.42706580696472
$ws.Range("B25").Value = 0.09235509630248373
$ws.Range("D25").Value = 0.0129191010734715
$ws.Range("E25").Value = 0.4984483718337742
$ws.Range("F25").Value = 0.818449577447808
$ws.Range("G25").Value = 0.6744398819355979
$ws.Range("H25").Value = 0.702885259457247
$ws.Range("I25").Value = 0.4313424922703462
$ws.Range("K25").Value = 1.17504717920059
